$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25-30 down to 26-31)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly price entry
$ws.Cells.Item(25, 1).Value = 1
$ws.Cells.Item(25, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(25, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(25, 4).Value = 45154
$ws.Cells.Item(25, 5).Value = 15
$ws.Cells.Item(25, 6).Value = 100114007
$ws.Cells.Item(25, 7).Value = "Jengibre"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 250
$ws.Cells.Item(25, 11).Value = 17000
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = 17500
$ws.Cells.Item(25, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 1346
$ws.Cells.Item(25, 17).Value = 13
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Ensure the date column keeps the date number format used by the rest of the column
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
